$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Q2: add new shared string text
$ws.Range("Q2").Value = "Wilayah_05_202206_31.zip"

# L3: value change 1121 -> 1177
$ws.Range("L3").Value = 1177

# Column Q width change
$ws.Columns("Q").ColumnWidth = 24.7109375

# Selection / view changes
$ws.Range("L3").Select()
$excel.ActiveWindow.ScrollColumn = 6
